$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.945.41'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '3.531.38'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.12%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '603.78'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.82%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '196.84'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E6').Value = '  +5.69%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.626'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -3.03%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.657'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E10').Value = '  +0.22%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '54.08'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.44%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0000303'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E12').Value = '  -2.54%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '4.086.09'
$ws.Range('E14').Value = '  -0.68%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '601.54'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E15').Value = '  -4.18%  '
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '70.086.71'
$ws.Range('E17').Value = '  -0.08%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '12.68'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '3.536.71'
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('E20').Value = '  +0.55%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.21%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '18.24'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E22').Value = '  +3.64%  '
$ws.Range('E23').Value = '  +7.19%  '
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('E26').Value = '  +2.40%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.96'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('E28').Value = '  +1.27%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '33.63'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.08%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.51'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E30').Value = '  +23.22%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.12'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E31').Value = '  +0.40%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '12.74'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('E33').Value = '  +1.45%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '63.51'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('E35').Value = '  +6.15%  '
$ws.Range('D36').Value = '3.741.30'
$ws.Range('E36').Value = '  +4.59%  '
$ws.Range('E37').Value = '  -5.32%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -2.01%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '36.87'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E41').Value = '  -1.18%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '494.29'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E42').Value = '  -6.93%  '
$ws.Range('E43').Value = '  -0.85%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0457'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.94%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.NumberFormat = 'General'
$c.ClearFormats()
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('E49').Value = '  -5.66%  '
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('E51').Value = '  +11.97%  '
